$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("N2").Value = "2018-12-31 00:00:00"

$ws.Range("O2").Value = 18361214792.96
$ws.Range("P2").Value = 929373239.6799999
$ws.Range("Q2").Value = 2641281831.96
$ws.Range("R2").Value = -37.5558970092
$ws.Range("S2").Value = 1173353768.3
$ws.Range("T2").Value = -29.5095030483
$ws.Range("U2").Value = 1227197755.46
$ws.Range("V2").Value = 11.4390860527
$ws.Range("W2").Value = 8852893214.139999
$ws.Range("X2").Value = 925230420.55
$ws.Range("Y2").Value = 28.7820161774
$ws.Range("Z2").Value = 1190971523.51
$ws.Range("AA2").Value = 82.8877653166
$ws.Range("AB2").Value = 9508321578.82
$ws.Range("AC2").Value = -9.893257466
$ws.Range("AD2").Value = -8.898426908899999
$ws.Range("AE2").Value = -7.8051832472
$ws.Range("AF2").Value = 102.7245160171
$ws.Range("AG2").Value = 48.2151824591
